# "Done, it is playable" commit
#  - slide 2: expand the "Button 1" caption textbox and rewrite its body copy
#  - slide 3: rewrite / restructure the "HOW TO PLAY" body copy, resize shape

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 2: "Button 1" textbox nested inside the "Gruppieren 6" group
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$grp = $s2.Shapes.Item(4)             # "Gruppieren 6"
$btn1 = $grp.GroupItems.Item(2)       # "Textfeld 4" -> Button 1 caption

# resize: cy 415498 -> 738664 EMU  (points = EMU / 12700)
$btn1.Height = 738664 / 12700

$btn1tr = $btn1.TextFrame.TextRange
$btn1tr.Text = "Button 1`rButton 1 does nothing, but if you are stressed, you are free to press it"
$btn1tr.Characters(1, 8).Font.Bold = $true

# ---------------------------------------------------------------------------
# Slide 3: "HOW TO PLAY" body textbox ("Textfeld 4")
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$body = $s3.Shapes.Item(2)            # "Textfeld 4"

# resize: cy 6247864 -> 5324535 EMU
$body.Height = 5324535 / 12700

$bodytr = $body.TextFrame.TextRange

$TAB = [string][char]9

$lines = @(
  "PLAYER SELECTION",
  "This game is designed as a single-player experience, so only one player can play at a time. There are no additional player settings. If you have a friend with you, you can simply take turns and see who achieves the higher score.",
  "@@P3@@",
  "OPTION SELECTION",
  "At the start of each game, after a brief introduction, you will arrive at the main menu. Here, you can choose from five different game modes:",
  ($TAB + "- STILL"),
  ($TAB + "- PULSE"),
  ($TAB + "- PORTAL"),
  ($TAB + "- ORBIT"),
  ($TAB + "- ECHO"),
  "@@P11@@",
  "A small arrow on the left side of the screen indicates the currently selected mode. The first time the game is launched, the default mode is STILL.",
  "You can navigate between the modes using the joystick. Push the joystick downward to move the arrow down. To continue moving in the same direction, release the joystick and then push it again in the desired direction (up or down).",
  "@@P14@@",
  "On the right side of the screen, you can see your current high score for the selected mode. A mode is considered completed when you reach a score of 16.",
  "@@P16@@",
  "To select a mode, press Button 4. You will immediately enter the level. If you decide you don't like the mode, you can return to the menu and select a different one."
)

$bodytr.Text = [string]::Join("`r", $lines)

# Clean up the placeholder-marked blank paragraphs (3, 11, 14, 16) by deleting
# the marker text, leaving a true empty paragraph behind.
foreach ($marker in @("@@P3@@", "@@P11@@", "@@P14@@", "@@P16@@")) {
  $full = $bodytr.Text
  $idx = $full.IndexOf($marker)
  $bodytr.Characters($idx + 1, $marker.Length).Text = ""
}

# ---------------------------------------------------------------------------
# Fix up run-level formatting. Setting the whole TextRange.Text above made
# every paragraph inherit paragraph 1's rPr (sz=1400 b=1); restore sz=1200
# on every paragraph but the two section headers, then re-apply bold only
# where the new copy wants it.
# ---------------------------------------------------------------------------

$full = $bodytr.Text

function Get-Idx($text, $needle, $from) {
  return $text.IndexOf($needle, $from)
}

# sz=1200 for everything from "This game is designed..." through the end.
$bodyStart = Get-Idx $full "This game is designed" 0
$bodyLen = $full.Length - $bodyStart
$bodytr.Characters($bodyStart + 1, $bodyLen).Font.Size = 12

# Turn off the inherited bold for that whole span too; the specific bold
# substrings (section header + list items + emphasised words) get re-applied
# below, layered on top.
$bodytr.Characters($bodyStart + 1, $bodyLen).Font.Bold = $false

# "OPTION SELECTION" header stays bold / sz=1400
$idx = Get-Idx $full "OPTION SELECTION" 0
$bodytr.Characters($idx + 1, "OPTION SELECTION".Length).Font.Bold = $true
$bodytr.Characters($idx + 1, "OPTION SELECTION".Length).Font.Size = 14

# Bulleted mode list: "- STILL" / "- PULSE" / "- PORTAL" / "- ORBIT" / "- ECHO" (bold)
$idx = Get-Idx $full "- STILL" 0
$bodytr.Characters($idx + 1, "- STILL".Length).Font.Bold = $true

$idx = Get-Idx $full "$TAB- PULSE" 0
$bodytr.Characters($idx + 1, ("$TAB- PULSE").Length).Font.Bold = $true

$idx = Get-Idx $full "$TAB- PORTAL" 0
$bodytr.Characters($idx + 1, ("$TAB- PORTAL").Length).Font.Bold = $true

$idx = Get-Idx $full "$TAB- ORBIT" 0
$bodytr.Characters($idx + 1, ("$TAB- ORBIT").Length).Font.Bold = $true

$idx = Get-Idx $full "$TAB- ECHO" 0
$bodytr.Characters($idx + 1, ("$TAB- ECHO").Length).Font.Bold = $true

# "STILL" emphasised inside "...the default mode is STILL."
$idx = Get-Idx $full "is STILL." 0
$idx = $idx + "is ".Length
$bodytr.Characters($idx + 1, "STILL".Length).Font.Bold = $true

# "completed" emphasised inside "...A mode is considered completed when..."
$idx = Get-Idx $full "considered completed" 0
$idx = $idx + "considered ".Length
$bodytr.Characters($idx + 1, "completed".Length).Font.Bold = $true

# "Button 4" emphasised inside "To select a mode, press Button 4."
$idx = Get-Idx $full "press Button 4" 0
$idx = $idx + "press ".Length
$bodytr.Characters($idx + 1, "Button 4".Length).Font.Bold = $true

Write-Host "Done."
Write-Host $bodytr.Text
